# Apply "first_release_qoq" series update:
#  - rename header B1 shared string: "value" -> "first_release_value"
#  - refresh/extend the date (A) / value (B) series from 53 to 84 rows
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "first_release_value"

$ws.Cells.Item(2, 1).Value = 38398
$ws.Cells.Item(2, 2).Value = -1.019024872684525
$ws.Cells.Item(3, 1).Value = 38487
$ws.Cells.Item(3, 2).Value = 0.8876024303436765
$ws.Cells.Item(4, 1).Value = 38579
$ws.Cells.Item(4, 2).Value = 0.501281665261132
$ws.Cells.Item(5, 1).Value = 38671
$ws.Cells.Item(5, 2).Value = 0.783791364788172
$ws.Cells.Item(6, 1).Value = 38763
$ws.Cells.Item(6, 2).Value = 0.2322508451440228
$ws.Cells.Item(7, 1).Value = 38852
$ws.Cells.Item(7, 2).Value = 0.8665864397470244
$ws.Cells.Item(8, 1).Value = 38944
$ws.Cells.Item(8, 2).Value = 0.6793221546917749
$ws.Cells.Item(9, 1).Value = 39036
$ws.Cells.Item(9, 2).Value = -0.9228048723025069
$ws.Cells.Item(10, 1).Value = 39128
$ws.Cells.Item(10, 2).Value = 1.742611766026243
$ws.Cells.Item(11, 1).Value = 39217
$ws.Cells.Item(11, 2).Value = -0.8957516004554691
$ws.Cells.Item(12, 1).Value = 39309
$ws.Cells.Item(12, 2).Value = 0.893923595129948
$ws.Cells.Item(13, 1).Value = 39401
$ws.Cells.Item(13, 2).Value = -0.2567725410682868
$ws.Cells.Item(14, 1).Value = 39493
$ws.Cells.Item(14, 2).Value = 1.910893826230975
$ws.Cells.Item(15, 1).Value = 39583
$ws.Cells.Item(15, 2).Value = -0.9623800849065276
$ws.Cells.Item(16, 1).Value = 39675
$ws.Cells.Item(16, 2).Value = 1.320073940737448
$ws.Cells.Item(17, 1).Value = 39767
$ws.Cells.Item(17, 2).Value = -0.05810994121875979
$ws.Cells.Item(18, 1).Value = 39859
$ws.Cells.Item(18, 2).Value = -1.734840982136873
$ws.Cells.Item(19, 1).Value = 39948
$ws.Cells.Item(19, 2).Value = -1.4178810011895
$ws.Cells.Item(20, 1).Value = 40040
$ws.Cells.Item(20, 2).Value = 1.258487819346726
$ws.Cells.Item(21, 1).Value = 40132
$ws.Cells.Item(21, 2).Value = -1.470147873721189
$ws.Cells.Item(22, 1).Value = 40224
$ws.Cells.Item(22, 2).Value = 1.38192517089017
$ws.Cells.Item(23, 1).Value = 40313
$ws.Cells.Item(23, 2).Value = 1.938408417700344
$ws.Cells.Item(24, 1).Value = 40405
$ws.Cells.Item(24, 2).Value = 0.4
$ws.Cells.Item(25, 1).Value = 40497
$ws.Cells.Item(25, 2).Value = -0.163634158232469
$ws.Cells.Item(26, 1).Value = 40589
$ws.Cells.Item(26, 2).Value = 1.079830393426633
$ws.Cells.Item(27, 1).Value = 40678
$ws.Cells.Item(27, 2).Value = 0.2
$ws.Cells.Item(28, 1).Value = 40770
$ws.Cells.Item(28, 2).Value = 0.4
$ws.Cells.Item(29, 1).Value = 40862
$ws.Cells.Item(29, 2).Value = 0.2303676816657827
$ws.Cells.Item(30, 1).Value = 40954
$ws.Cells.Item(30, 2).Value = -0.3
$ws.Cells.Item(31, 1).Value = 41044
$ws.Cells.Item(31, 2).Value = -0.4045626487644824
$ws.Cells.Item(32, 1).Value = 41136
$ws.Cells.Item(32, 2).Value = 0
$ws.Cells.Item(33, 1).Value = 41228
$ws.Cells.Item(33, 2).Value = -0.01847206600469065
$ws.Cells.Item(34, 1).Value = 41320
$ws.Cells.Item(34, 2).Value = 0.009235986179263023
$ws.Cells.Item(35, 1).Value = 41409
$ws.Cells.Item(35, 2).Value = 0.3863474960573257
$ws.Cells.Item(36, 1).Value = 41501
$ws.Cells.Item(36, 2).Value = 0.6964197943645729
$ws.Cells.Item(37, 1).Value = 41593
$ws.Cells.Item(37, 2).Value = -0.3183250015702015
$ws.Cells.Item(38, 1).Value = 41685
$ws.Cells.Item(38, 2).Value = 1.852186157158073
$ws.Cells.Item(39, 1).Value = 41774
$ws.Cells.Item(39, 2).Value = -0.1726927221574073
$ws.Cells.Item(40, 1).Value = 41866
$ws.Cells.Item(40, 2).Value = -0.1729977607768376
$ws.Cells.Item(41, 1).Value = 41958
$ws.Cells.Item(41, 2).Value = 1.050807574684342
$ws.Cells.Item(42, 1).Value = 42050
$ws.Cells.Item(42, 2).Value = 0.5056288600178789
$ws.Cells.Item(43, 1).Value = 42139
$ws.Cells.Item(43, 2).Value = -0.2184803162966205
$ws.Cells.Item(44, 1).Value = 42231
$ws.Cells.Item(44, 2).Value = 0.7235341094351355
$ws.Cells.Item(45, 1).Value = 42323
$ws.Cells.Item(45, 2).Value = 0.8513884674671885
$ws.Cells.Item(46, 1).Value = 42415
$ws.Cells.Item(46, 2).Value = 0.8066734233961483
$ws.Cells.Item(47, 1).Value = 42505
$ws.Cells.Item(47, 2).Value = -0.1388227614901609
$ws.Cells.Item(48, 1).Value = 42597
$ws.Cells.Item(48, 2).Value = 0.4819278240608753
$ws.Cells.Item(49, 1).Value = 42689
$ws.Cells.Item(49, 2).Value = 0.6629265129002277
$ws.Cells.Item(50, 1).Value = 42781
$ws.Cells.Item(50, 2).Value = 0.1646390629436354
$ws.Cells.Item(51, 1).Value = 42870
$ws.Cells.Item(51, 2).Value = 1.164700738417963
$ws.Cells.Item(52, 1).Value = 42962
$ws.Cells.Item(52, 2).Value = 0.4317240674915439
$ws.Cells.Item(53, 1).Value = 43054
$ws.Cells.Item(53, 2).Value = 0.1359170431485039
$ws.Cells.Item(54, 1).Value = 43146
$ws.Cells.Item(54, 2).Value = 0.3857269132374052
$ws.Cells.Item(55, 1).Value = 43235
$ws.Cells.Item(55, 2).Value = 0.876336956515118
$ws.Cells.Item(56, 1).Value = 43327
$ws.Cells.Item(56, 2).Value = 0.8
$ws.Cells.Item(57, 1).Value = 43419
$ws.Cells.Item(57, 2).Value = 0.008724159582257585
$ws.Cells.Item(58, 1).Value = 43511
$ws.Cells.Item(58, 2).Value = 0.2
$ws.Cells.Item(59, 1).Value = 43600
$ws.Cells.Item(59, 2).Value = 0.4878538807911497
$ws.Cells.Item(60, 1).Value = 43692
$ws.Cells.Item(60, 2).Value = -0.4043302599539206
$ws.Cells.Item(61, 1).Value = 43784
$ws.Cells.Item(61, 2).Value = 0.7183553771707381
$ws.Cells.Item(62, 1).Value = 43876
$ws.Cells.Item(62, 2).Value = -1.5
$ws.Cells.Item(63, 1).Value = 43966
$ws.Cells.Item(63, 2).Value = -7.231044133207007
$ws.Cells.Item(64, 1).Value = 44058
$ws.Cells.Item(64, 2).Value = 4.729401638091318
$ws.Cells.Item(65, 1).Value = 44150
$ws.Cells.Item(65, 2).Value = -0.8905127363963885
$ws.Cells.Item(66, 1).Value = 44242
$ws.Cells.Item(66, 2).Value = -2.04269378128221
$ws.Cells.Item(67, 1).Value = 44331
$ws.Cells.Item(67, 2).Value = 2.093024636165651
$ws.Cells.Item(68, 1).Value = 44423
$ws.Cells.Item(68, 2).Value = 1.098535546956398
$ws.Cells.Item(69, 1).Value = 44515
$ws.Cells.Item(69, 2).Value = -0.4717175472572421
$ws.Cells.Item(70, 1).Value = 44607
$ws.Cells.Item(70, 2).Value = 0.9401304606753627
$ws.Cells.Item(71, 1).Value = 44696
$ws.Cells.Item(71, 2).Value = 0.6392725048137464
$ws.Cells.Item(72, 1).Value = 44788
$ws.Cells.Item(72, 2).Value = 0.06357296580725347
$ws.Cells.Item(73, 1).Value = 44880
$ws.Cells.Item(73, 2).Value = -1.012166871044968
$ws.Cells.Item(74, 1).Value = 44972
$ws.Cells.Item(74, 2).Value = -0.5703626997413522
$ws.Cells.Item(75, 1).Value = 45061
$ws.Cells.Item(75, 2).Value = 0.2394101325822788
$ws.Cells.Item(76, 1).Value = 45153
$ws.Cells.Item(76, 2).Value = -0.04072131480353391
$ws.Cells.Item(77, 1).Value = 45245
$ws.Cells.Item(77, 2).Value = -0.07958838003274593
$ws.Cells.Item(78, 1).Value = 45337
$ws.Cells.Item(78, 2).Value = 0.02912383308249389
$ws.Cells.Item(79, 1).Value = 45427
$ws.Cells.Item(79, 2).Value = -0.1311265493919933
$ws.Cells.Item(80, 1).Value = 45519
$ws.Cells.Item(80, 2).Value = -0.03907468377752821
$ws.Cells.Item(81, 1).Value = 45611
$ws.Cells.Item(81, 2).Value = 0.1771324545010202
$ws.Cells.Item(82, 1).Value = 45703
$ws.Cells.Item(82, 2).Value = 0.4946531409412387
$ws.Cells.Item(83, 1).Value = 45792
$ws.Cells.Item(83, 2).Value = 0.202428137729683
$ws.Cells.Item(84, 1).Value = 45884
$ws.Cells.Item(84, 2).Value = 0.208573386070384

# Rows 54-84 are brand new; give column A the same date number-format
# style already used by the rest of the series (copied from A2).
$ws.Range("A2").Copy()
$ws.Range("A54:A84").PasteSpecial(-4122)
$excel.CutCopyMode = $false
